# Atualiza dados da BIBI - faturamento anual (linha 2025)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = 3622402.33
$ws.Range("C9").Value = 570775.24
$ws.Range("D9").Value = 4193177.57
$ws.Range("E9").Value = 13.61199783390046
$ws.Range("F9").Value = 86.38800216609953
$ws.Range("G9").Value = -44.83741907055273
$ws.Range("H9").Value = -34.5844209482003
$ws.Range("I9").Value = 36570
$ws.Range("J9").Value = 1562
$ws.Range("K9").Value = 38132
$ws.Range("L9").Value = 26341
$ws.Range("M9").Value = 159.1882453209825
$ws.Range("N9").Value = 8.681057383740676
